$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.208.47"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.860.18"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.80"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06580"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.89"
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07986"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.68"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "1.860.64"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.114"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.81"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").Value = "30.180.63"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.60"
$ws.Range("E18").Value = "  +6.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007660"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "2.103.85"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -4.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.12"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.195"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.98"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.947"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.374"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09912"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.336"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.465"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.039"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04716"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.124"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7022"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.343"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.59"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.937"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8399"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4154"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9977"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.61"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.238"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.071"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "935.25"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.09"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05661"
$ws.Range("E51").Value = "  +0.49%  "
